$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$paragraph.Range.InsertXML($pkg)
}

# "TABLA DE CONTENIDO" heading paragraph: drop the red-color rPr from both
# the paragraph mark (pPr) and the run, keeping every other attribute as-is.
$pTabla = $d.Paragraphs.Item(32)
$xmlTabla = '<w:p w14:paraId="6380D1DA" w14:textId="71A977E2" w:rsidR="007675FD" w:rsidRPr="000A5BCA" w:rsidRDefault="003E515B" w:rsidP="00A23292">' + '<w:pPr><w:jc w:val="center"/></w:pPr>' + '<w:r w:rsidRPr="000A5BCA"><w:lastRenderedPageBreak/><w:t>TABLA DE CONTENIDO</w:t></w:r>' + '</w:p>'
Set-ParagraphXml $pTabla $xmlTabla

# "1.dd" entry: the pPr held only the color rPr, so it disappears entirely.
$p1 = $d.Paragraphs.Item(33)
$xml1 = '<w:p w14:paraId="4DCD9CDA" w14:textId="37CA560E" w:rsidR="003E515B" w:rsidRPr="000A5BCA" w:rsidRDefault="003E515B" w:rsidP="003E515B">' + '<w:r w:rsidRPr="000A5BCA"><w:t>1.dd</w:t></w:r>' + '</w:p>'
Set-ParagraphXml $p1 $xml1

# "2.ss" entry: same pattern.
$p2 = $d.Paragraphs.Item(34)
$xml2 = '<w:p w14:paraId="1D6FDED8" w14:textId="116BE459" w:rsidR="003E515B" w:rsidRPr="000A5BCA" w:rsidRDefault="003E515B" w:rsidP="003E515B">' + '<w:r w:rsidRPr="000A5BCA"><w:t>2.ss</w:t></w:r>' + '</w:p>'
Set-ParagraphXml $p2 $xml2

# "3.mm" entry: same pattern.
$p3 = $d.Paragraphs.Item(35)
$xml3 = '<w:p w14:paraId="52AD0F29" w14:textId="4431A6A7" w:rsidR="003E515B" w:rsidRPr="000A5BCA" w:rsidRDefault="003E515B" w:rsidP="003E515B">' + '<w:r w:rsidRPr="000A5BCA"><w:t>3.mm</w:t></w:r>' + '</w:p>'
Set-ParagraphXml $p3 $xml3

# Trailing empty paragraph after the table of contents: its pPr held only
# the color rPr too, so it becomes a fully empty <w:p/>.
$p4 = $d.Paragraphs.Item(36)
$xml4 = '<w:p w14:paraId="75BD937F" w14:textId="77777777" w:rsidR="003E515B" w:rsidRPr="000A5BCA" w:rsidRDefault="003E515B" w:rsidP="003E515B"/>'
Set-ParagraphXml $p4 $xml4

Write-Host "done"
